$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count
$firstRow = $used.Row
$firstCol = $used.Column

for ($r = 0; $r -lt $rowCount; $r++) {
    for ($c = 0; $c -lt $colCount; $c++) {
        $cell = $ws.Cells.Item($firstRow + $r, $firstCol + $c)
        $val = $cell.Value2
        if ($val -is [string]) {
            $newVal = $val.Replace("D51", "D55").Replace("D64", "D69").Replace("D80", "D86").Replace("S30", "S31")
            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
